# Formed the consolidated report
# Update the "Absent" (column H) values on Sheet1 to reflect the
# consolidated attendance report: rows that previously had no value
# (inlineStr placeholder) now get an explicit 0, and rows whose Absent
# flag needed correcting are set to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H13").Value = 0
